$wb = $excel.ActiveWorkbook

# ---- Worksheet 1 ----
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2,1).Value = "Última actualización: 08:28:52"
$ws.Cells.Item(3,1).Value = "Total filas: 71"

$ws.Rows.Item(40).Insert()
$ws.Cells.Item(40,1).Value = "08:28:52"
$ws.Cells.Item(40,2).Value = "08:32"
$ws.Cells.Item(40,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(40,4).Value = 4
$ws.Cells.Item(40,5).Value = "LP1912"

$ws.Rows.Item(50).Insert()
$ws.Cells.Item(50,1).Value = "08:28:52"
$ws.Cells.Item(50,2).Value = "08:57"
$ws.Cells.Item(50,3).Value = "215A_EL PATO"
$ws.Cells.Item(50,4).Value = 29
$ws.Cells.Item(50,5).Value = "LP1912"

$ws.Rows.Item(52).Insert()
$ws.Cells.Item(52,1).Value = "08:28:52"
$ws.Cells.Item(52,2).Value = "09:04"
$ws.Cells.Item(52,3).Value = "10_OLMOS"
$ws.Cells.Item(52,4).Value = 36
$ws.Cells.Item(52,5).Value = "LP1912"

$ws.Rows.Item(58).Insert()
$ws.Cells.Item(58,1).Value = "08:28:52"
$ws.Cells.Item(58,2).Value = "09:16"
$ws.Cells.Item(58,3).Value = "27_EL RETIRO"
$ws.Cells.Item(58,4).Value = 48
$ws.Cells.Item(58,5).Value = "LP1912"

$ws.Rows.Item(63).Insert()
$ws.Cells.Item(63,1).Value = "08:28:52"
$ws.Cells.Item(63,2).Value = "09:28"
$ws.Cells.Item(63,3).Value = "10_OLMOS"
$ws.Cells.Item(63,4).Value = 60
$ws.Cells.Item(63,5).Value = "LP1912"

$ws.Rows.Item(66).Insert()
$ws.Cells.Item(66,1).Value = "08:28:52"
$ws.Cells.Item(66,2).Value = "09:32"
$ws.Cells.Item(66,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(66,4).Value = 64
$ws.Cells.Item(66,5).Value = "LP1912"

$ws.Rows.Item(72).Insert()
$ws.Cells.Item(72,1).Value = "08:28:52"
$ws.Cells.Item(72,2).Value = "09:58"
$ws.Cells.Item(72,3).Value = "215C_EL PATO"
$ws.Cells.Item(72,4).Value = 90
$ws.Cells.Item(72,5).Value = "LP1912"

$ws.Rows.Item(74).Insert()
$ws.Cells.Item(74,1).Value = "08:28:52"
$ws.Cells.Item(74,2).Value = "10:05"
$ws.Cells.Item(74,3).Value = "14_ABASTO"
$ws.Cells.Item(74,4).Value = 97
$ws.Cells.Item(74,5).Value = "LP1912"

$ws.Rows.Item(76).Insert()
$ws.Cells.Item(76,1).Value = "08:28:52"
$ws.Cells.Item(76,2).Value = "10:13"
$ws.Cells.Item(76,3).Value = "17X38_ROMERO"
$ws.Cells.Item(76,4).Value = 105
$ws.Cells.Item(76,5).Value = "LP1912"

# ---- Worksheet 2 ----
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2,1).Value = "Última actualización: 08:28:52"
$ws.Cells.Item(3,1).Value = "Total filas: 16"

$ws.Rows.Item(18).Insert()
$ws.Cells.Item(18,1).Value = "08:28:52"
$ws.Cells.Item(18,2).Value = "08:57"
$ws.Cells.Item(18,3).Value = "215A_EL PATO"
$ws.Cells.Item(18,4).Value = 29
$ws.Cells.Item(18,5).Value = "LP1912"

$ws.Rows.Item(20).Insert()
$ws.Cells.Item(20,1).Value = "08:28:52"
$ws.Cells.Item(20,2).Value = "09:58"
$ws.Cells.Item(20,3).Value = "215C_EL PATO"
$ws.Cells.Item(20,4).Value = 90
$ws.Cells.Item(20,5).Value = "LP1912"

# ---- Worksheet 3 ----
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2,1).Value = "Última actualización: 08:28:52"
$ws.Cells.Item(3,1).Value = "Total filas: 9"

$ws.Rows.Item(12).Insert()
$ws.Cells.Item(12,1).Value = "08:28:52"
$ws.Cells.Item(12,2).Value = "09:20"
$ws.Cells.Item(12,3).Value = "215A_LA PLATA"
$ws.Cells.Item(12,4).Value = 52
$ws.Cells.Item(12,5).Value = "L6173"

$ws.Rows.Item(14).Insert()
$ws.Cells.Item(14,1).Value = "08:28:52"
$ws.Cells.Item(14,2).Value = "10:12"
$ws.Cells.Item(14,3).Value = "215C_LA PLATA"
$ws.Cells.Item(14,4).Value = 104
$ws.Cells.Item(14,5).Value = "L6203"
